$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("What is the maximum number of tracks that can be specified within a single ODF file?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."),
    @("What is the maximum number of tracks that can be specified within a single ODF file?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."),
    @("How many tracks can you set up in one ODF?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."),
    @("What’s the track limit for an ODF?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most. It's recommended to use OIF files instead, which allows you to save only the specific depth range of interest, making data management more efficient."),
    @("In one ODF configuration, how many tracks can be defined?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."),
    @("Is there a maximum number of tracks you can create in one ODF?", "llama3.2:latest", "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most.")
)

$startRow = 116
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
